$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.675.61"
$ws.Range("E2").Value = "  +3.97%  "
$ws.Range("D3").Value = "3.079.08"
$ws.Range("E3").Value = "  +3.68%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'551.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.42%  "
$ws.Range("D6").Value = "'138.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.59%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.072.49"
$ws.Range("E8").Value = "  +3.57%  "
$ws.Range("D9").Value = "'0.500"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.40%  "
$ws.Range("D10").Value = "'0.150"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("D11").Value = "'6.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("D12").Value = "'0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.85%  "
$ws.Range("D13").Value = "'0.0000226"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.68%  "
$ws.Range("D14").Value = "'34.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.70%  "
$ws.Range("D15").Value = "3.571.79"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").Value = "63.610.18"
$ws.Range("E16").Value = "  +3.69%  "
$ws.Range("D17").Value = "3.082.26"
$ws.Range("E17").Value = "  +3.71%  "
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").Value = "'6.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.78%  "
$ws.Range("D20").Value = "'486.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.88%  "
$ws.Range("D21").Value = "'13.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.67%  "
$ws.Range("D22").Value = "'0.685"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.64%  "
$ws.Range("D23").Value = "'7.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.23%  "
$ws.Range("D24").Value = "'81.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.80%  "
$ws.Range("E25").Value = "  +5.86%  "
$ws.Range("D27").Value = "'2.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.60%  "
$ws.Range("D28").Value = "'8.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.87%  "
$ws.Range("D29").Value = "'1.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.89%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "'26.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.97%  "
$ws.Range("D32").Value = "'1.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("D33").Value = "'5.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.64%  "
$ws.Range("D34").Value = "'2.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.72%  "
$ws.Range("D35").Value = "'55.78"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("D36").Value = "'5.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.21%  "
$ws.Range("D37").Value = "'471.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.10%  "
$ws.Range("D38").Value = "3.186.11"
$ws.Range("E38").Value = "  +2.07%  "
$ws.Range("D39").Value = "'0.0818"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.90%  "
$ws.Range("D40").Value = "'0.0396"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.23%  "
$ws.Range("E41").Value = "  +5.07%  "
$ws.Range("D42").Value = "'8.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.96%  "
$ws.Range("D43").Value = "'2.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.53%  "
$ws.Range("D44").Value = "'27.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +13.02%  "
$ws.Range("D45").Value = "'0.252"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.41%  "
$ws.Range("D47").Value = "'2.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.09%  "
$ws.Range("E48").Value = "  +3.42%  "
$ws.Range("D49").Value = "0.0₃0514"
$ws.Range("E49").Value = "  +2.83%  "
$ws.Range("D50").Value = "'116.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.98%  "
$ws.Range("D51").Value = "'2.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.77%  "
